$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-5) are being re-sorted ascending by Fecha (column D).
# Row 2 <-> Row 4 swap, Row 3 <-> Row 5 swap for columns D, L, M, N, O, P, S.

# --- Save original values for rows 2 and 3 before overwriting ---
$D2 = $ws.Range("D2").Value2
$L2 = $ws.Range("L2").Value2
$M2 = $ws.Range("M2").Value2
$N2 = $ws.Range("N2").Value2
$O2 = $ws.Range("O2").Value2
$P2 = $ws.Range("P2").Value2
$S2 = $ws.Range("S2").Value2

$D3 = $ws.Range("D3").Value2
$L3 = $ws.Range("L3").Value2
$M3 = $ws.Range("M3").Value2
$N3 = $ws.Range("N3").Value2
$O3 = $ws.Range("O3").Value2
$P3 = $ws.Range("P3").Value2
$S3 = $ws.Range("S3").Value2

# --- Save original values for rows 4 and 5 ---
$D4 = $ws.Range("D4").Value2
$L4 = $ws.Range("L4").Value2
$M4 = $ws.Range("M4").Value2
$N4 = $ws.Range("N4").Value2
$O4 = $ws.Range("O4").Value2
$P4 = $ws.Range("P4").Value2
$S4 = $ws.Range("S4").Value2

$D5 = $ws.Range("D5").Value2
$L5 = $ws.Range("L5").Value2
$M5 = $ws.Range("M5").Value2
$N5 = $ws.Range("N5").Value2
$O5 = $ws.Range("O5").Value2
$P5 = $ws.Range("P5").Value2
$S5 = $ws.Range("S5").Value2

# --- Write row 4's original data into row 2 ---
$ws.Range("D2").Value2 = $D4
$ws.Range("L2").Value2 = $L4
$ws.Range("M2").Value2 = $M4
$ws.Range("N2").Value2 = $N4
$ws.Range("O2").Value2 = $O4
$ws.Range("P2").Value2 = $P4
$ws.Range("S2").Value2 = $S4

# --- Write row 5's original data into row 3 ---
$ws.Range("D3").Value2 = $D5
$ws.Range("L3").Value2 = $L5
$ws.Range("M3").Value2 = $M5
$ws.Range("N3").Value2 = $N5
$ws.Range("O3").Value2 = $O5
$ws.Range("P3").Value2 = $P5
$ws.Range("S3").Value2 = $S5

# --- Write row 2's original data into row 4 ---
$ws.Range("D4").Value2 = $D2
$ws.Range("L4").Value2 = $L2
$ws.Range("M4").Value2 = $M2
$ws.Range("N4").Value2 = $N2
$ws.Range("O4").Value2 = $O2
$ws.Range("P4").Value2 = $P2
$ws.Range("S4").Value2 = $S2

# --- Write row 3's original data into row 5 ---
$ws.Range("D5").Value2 = $D3
$ws.Range("L5").Value2 = $L3
$ws.Range("M5").Value2 = $M3
$ws.Range("N5").Value2 = $N3
$ws.Range("O5").Value2 = $O3
$ws.Range("P5").Value2 = $P3
$ws.Range("S5").Value2 = $S3
